$d = $word.ActiveDocument

$replacements = @(
    @("86÷4=21, 2", "44÷7=6, 2"),
    @("55÷8=6, 7", "98÷7=14, 0"),
    @("87÷9=9, 6", "28÷8=3, 4"),
    @("55÷6=9, 1", "73÷3=24, 1"),
    @("53÷7=7, 4", "94÷4=23, 2"),
    @("47÷9=5, 2", "10÷5=2, 0"),
    @("95÷7=13, 4", "37÷3=12, 1"),
    @("22÷4=5, 2", "22÷5=4, 2"),
    @("73÷2=36, 1", "44÷9=4, 8"),
    @("65÷5=13, 0", "88÷9=9, 7"),
    @("52÷7=7, 3", "90÷2=45, 0"),
    @("17÷7=2, 3", "24÷3=8, 0"),
    @("85÷8=10, 5", "88÷6=14, 4"),
    @("80÷2=40, 0", "48÷7=6, 6"),
    @("40÷4=10, 0", "89÷7=12, 5"),
    @("29÷2=14, 1", "43÷6=7, 1"),
    @("15÷4=3, 3", "87÷4=21, 3"),
    @("98÷3=32, 2", "14÷7=2, 0"),
    @("24÷7=3, 3", "15÷3=5, 0"),
    @("12÷8=1, 4", "57÷2=28, 1"),
    @("11÷3=3, 2", "74÷2=37, 0"),
    @("77÷9=8, 5", "22÷5=4, 2"),
    @("31÷2=15, 1", "32÷9=3, 5"),
    @("36÷7=5, 1", "97÷9=10, 7"),
    @("83÷7=11, 6", "13÷5=2, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
